$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.6383931775788736,
    -23.10748969621655,
    0.4106974857706616,
    -0.7814325448240136,
    0.8293808000575305,
    0.2146654303905574,
    14.31124727300474,
    0.09521762600208575,
    0.1344493681223383,
    0.114833497062212,
    0.2718993657310753,
    0.4633200086231518,
    -0.08482046726337922,
    0.4830445005948765,
    35.07734920122535,
    54.57936239911656
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
